$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (FAPs -> FAPs) specificity values (TPM recompute)
$ws.Range("I2").Value = 0.02394963654761903
$ws.Range("J2").Value = 0.02394963654761903
$ws.Range("S2").Value = 0.02394963654761903
$ws.Range("T2").Value = 0.02394963654761903

# Row 3 now represents the Resolving-Mac -> FAPs edge (the old MuSCs row is
# dropped below), with new TPM-derived values
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 34.93245566666667
$ws.Range("H3").Value = 104.797367
$ws.Range("I3").Value = 0.9760503634523809
$ws.Range("J3").Value = 0.9760503634523809
$ws.Range("Q3").Value = 9.593407146521557
$ws.Range("R3").Value = 86.34066431869401
$ws.Range("S3").Value = 0.9760503634523809
$ws.Range("T3").Value = 0.9760503634523809

# Remove the now-obsolete row 4 (Resolving-Mac -> FAPs); its updated data
# now lives in row 3, and the "MuSCs" shared string becomes unused.
$ws.Rows.Item(4).Delete()
